$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Review orders: "Test order N" -> "Work order N" (dropdown -> checkbox display text)
$ws.Range("C3").Value = "Work order 1"
$ws.Range("C4").Value = "Work order 2"
$ws.Range("C5").Value = "Work order 3"
$ws.Range("C6").Value = "Work order 4"

# partId column now holds simple sequential numbers instead of raw part codes
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 3

# Those cells are no longer bold now that they reflect simple checkbox-linked order numbers
$ws.Range("B3:B6").Font.Bold = $false

# Widen column C to fit the new "Work order N" / checkbox labels
$ws.Columns("C").ColumnWidth = 23.25

# Move the active selection as it was left after the edit
[void]$ws.Range("C10").Select()
